# Commit: "changing document, table attributes to lowerCamelCase"
#
# The ObjTables header rows embedded as plain text in row 1 (and, on the
# first sheet, also row 2) use PascalCase attribute names
# (ObjTablesVersion=, Type=, Id=). They must become lowerCamelCase
# (objTablesVersion=, type=, id=) while keeping everything else in the
# string identical.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "!!FirstUnambiguousModel" -> content says FirstAmbiguousModel ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$ws1.Range("A2").Value = "!!ObjTables type='Data' id='FirstAmbiguousModel'"

# --- Sheet 2: "!!SecondUnambiguousModel" -> content says SecondAmbiguousModel ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Value = "!!ObjTables type='Data' id='SecondAmbiguousModel'"
# Touch row 2 so it becomes part of the sheet's used range (dimension
# A1:A1 -> A1:A2) without giving it any real content/value.
$ws2.Range("A2").Style = "Normal"

# --- Sheet 3: "!!TestModel" ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A1").Value = "!!ObjTables type='Data' id='TestModel'"
$ws3.Range("A2").Style = "Normal"

# --- Sheet 4: "!!TestModels3" -> content says TestModel3 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A1").Value = "!!ObjTables type='Data' id='TestModel3'"
$ws4.Range("A2").Style = "Normal"
